$d = $word.ActiveDocument
$d.Content.Find.Execute("2024-01-08 Monday", $true, $false, $false, $false, $false, $true, 1, $false, "2024-01-09 Tuesday", 2) | Out-Null
$t = $d.Tables.Item(1)
$t.Cell(1, 1).Range.Text = "86×71="
$t.Cell(1, 2).Range.Text = "94×53="
$t.Cell(1, 3).Range.Text = "55×90="
$t.Cell(1, 4).Range.Text = "98×50="
$t.Cell(1, 5).Range.Text = "85×96="
$t.Cell(5, 1).Range.Text = "33×63="
$t.Cell(5, 2).Range.Text = "32×35="
$t.Cell(5, 3).Range.Text = "43×85="
$t.Cell(5, 4).Range.Text = "22×81="
$t.Cell(5, 5).Range.Text = "31×72="
$t.Cell(10, 1).Range.Text = "39×67="
$t.Cell(10, 2).Range.Text = "71×37="
$t.Cell(10, 3).Range.Text = "12×71="
$t.Cell(10, 4).Range.Text = "43×87="
$t.Cell(10, 5).Range.Text = "94×37="
$t.Cell(15, 1).Range.Text = "46×49="
$t.Cell(15, 2).Range.Text = "40×74="
$t.Cell(15, 3).Range.Text = "57×28="
$t.Cell(15, 4).Range.Text = "41×23="
$t.Cell(15, 5).Range.Text = "58×69="
$t.Cell(20, 1).Range.Text = "86×19="
$t.Cell(20, 2).Range.Text = "62×45="
$t.Cell(20, 3).Range.Text = "43×66="
$t.Cell(20, 4).Range.Text = "46×34="
$t.Cell(20, 5).Range.Text = "65×45="
